# Update TAKE-OFF sheet: move the CO2/SOx emission values from rows 44/45
# up to rows 40/41 (NOx / CO emissions), zeroing out the old rows 44/45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TAKE-OFF")

$ws.Range("C40").Value = 30111.785835342387
$ws.Range("C41").Value = 11903.773654072735
$ws.Range("C44").Value = 0.0
$ws.Range("C45").Value = 0.0
